# Correcting market share tab for updated scenario 3s
#
# The "New Product A / MDA" row (row 2) on the MarketShare sheet should
# carry the 100% market share for years 2026-2040 (columns L:Z), while the
# "Old Product B (SOC) / MDA" row (row 3) keeps the 100% market share only
# for years 2018-2025 (columns D:K) - i.e. the years 2026-2040 values move
# from row 3 to row 2.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("MarketShare")

# Move the "1" values for 2026-2040 (columns L:Z) from row 3 to row 2.
$ws2.Range("L2:Z2").Value = 1
$ws2.Range("L3:Z3").ClearContents()

# Make the MarketShare tab the active / displayed sheet, with the
# selection sitting on the newly-populated range.
$ws2.Activate()
$ws2.Range("L2:Z2").Select()
